# Auto-generated Excel COM-interop script
# Applies refreshed market-data values (currentAveragePrice*, LevePrice*, LeveProfit*)
# to the Cactuar_Profits workbook sheets, per the scheduled-runner data update.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1009.6
$ws.Range("I19").Value = 950
$ws.Range("J19").Value = 1049.3334
$ws.Range("K19").Value = 950
$ws.Range("L19").Value = 1049.3334
$ws.Range("M19").Value = -775
$ws.Range("N19").Value = -1399.3334
$ws.Range("H33").Value = 849.5714
$ws.Range("I33").Value = 849.5714
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 849.5714
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -620.5714
$ws.Range("H38").Value = 4491.9287
$ws.Range("J38").Value = 8678
$ws.Range("L38").Value = 26034
$ws.Range("N38").Value = -26778
$ws.Range("H43").Value = 2568100.8
$ws.Range("J43").Value = 4747.25
$ws.Range("L43").Value = 4747.25
$ws.Range("N43").Value = -4885.25
$ws.Range("H70").Value = 2385.0908
$ws.Range("J70").Value = 1842
$ws.Range("L70").Value = 5526
$ws.Range("N70").Value = -6066
$ws.Range("H73").Value = 2385.0908
$ws.Range("J73").Value = 1842
$ws.Range("L73").Value = 5526
$ws.Range("N73").Value = -7398
$ws.Range("H116").Value = 48195630
$ws.Range("J116").Value = 71435864
$ws.Range("L116").Value = 71435864
$ws.Range("N116").Value = -71442748
$ws.Range("H129").Value = 1409.4814
$ws.Range("I129").Value = 724.3333
$ws.Range("J129").Value = 2265.9167
$ws.Range("K129").Value = 2172.9999
$ws.Range("L129").Value = 6797.750100000001
$ws.Range("M129").Value = 2827.0001
$ws.Range("N129").Value = -16797.7501
$ws.Range("H132").Value = 14398.269
$ws.Range("I132").Value = 2885.862
$ws.Range("J132").Value = 23184.053
$ws.Range("K132").Value = 8657.585999999999
$ws.Range("L132").Value = 69552.159
$ws.Range("M132").Value = -6127.585999999999
$ws.Range("N132").Value = -74612.159
$ws.Range("H141").Value = 5781.5454
$ws.Range("I141").Value = 7832.3335
$ws.Range("K141").Value = 23497.0005
$ws.Range("M141").Value = -18317.0005
$ws.Range("N33").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1250165.9
$ws.Range("I2").Value = 1590243.5
$ws.Range("J2").Value = 3214.6667
$ws.Range("K2").Value = 1590243.5
$ws.Range("L2").Value = 3214.6667
$ws.Range("M2").Value = -1590130.5
$ws.Range("N2").Value = -3440.6667
$ws.Range("H32").Value = 19326.244
$ws.Range("I32").Value = 21487.621
$ws.Range("K32").Value = 21487.621
$ws.Range("M32").Value = -21200.621
$ws.Range("H61").Value = 15867.846
$ws.Range("I61").Value = 19365.334
$ws.Range("K61").Value = 19365.334
$ws.Range("M61").Value = -19153.334
$ws.Range("H74").Value = 1012.52
$ws.Range("I74").Value = 811.3570999999999
$ws.Range("J74").Value = 1268.5454
$ws.Range("K74").Value = 811.3570999999999
$ws.Range("L74").Value = 1268.5454
$ws.Range("M74").Value = 62.64290000000005
$ws.Range("N74").Value = -3016.5454
$ws.Range("H77").Value = 1012.52
$ws.Range("I77").Value = 811.3570999999999
$ws.Range("J77").Value = 1268.5454
$ws.Range("K77").Value = 4056.7855
$ws.Range("L77").Value = 6342.727
$ws.Range("M77").Value = 311.2145
$ws.Range("N77").Value = -15078.727
$ws.Range("H101").Value = 99475
$ws.Range("J101").Value = 99475
$ws.Range("L101").Value = 99475
$ws.Range("H116").Value = 1250165.9
$ws.Range("I116").Value = 1590243.5
$ws.Range("J116").Value = 3214.6667
$ws.Range("K116").Value = 1590243.5
$ws.Range("L116").Value = 3214.6667
$ws.Range("M116").Value = -1587949.5
$ws.Range("N116").Value = -7802.6667
$ws.Range("H132").Value = 16010.293
$ws.Range("I132").Value = 17894.766
$ws.Range("J132").Value = 6857.143
$ws.Range("K132").Value = 53684.298
$ws.Range("L132").Value = 20571.429
$ws.Range("M132").Value = -51154.298
$ws.Range("N132").Value = -25631.429
$ws.Range("H136").Value = 15867.846
$ws.Range("I136").Value = 19365.334
$ws.Range("K136").Value = 58096.00199999999
$ws.Range("M136").Value = -55546.00199999999
$ws.Range("N101").Value = -105965

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1250165.9
$ws.Range("I3").Value = 1590243.5
$ws.Range("J3").Value = 3214.6667
$ws.Range("K3").Value = 1590243.5
$ws.Range("L3").Value = 3214.6667
$ws.Range("M3").Value = -1590129.5
$ws.Range("N3").Value = -3442.6667
$ws.Range("H105").Value = 55557350
$ws.Range("I105").Value = 90910776
$ws.Range("J105").Value = 1973.1428
$ws.Range("K105").Value = 90910776
$ws.Range("L105").Value = 1973.1428
$ws.Range("M105").Value = -90909029
$ws.Range("N105").Value = -5467.1428
$ws.Range("H134").Value = 3591.95
$ws.Range("I134").Value = 1488.75
$ws.Range("K134").Value = 4466.25
$ws.Range("M134").Value = -1931.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 26897
$ws.Range("J98").Value = 26897
$ws.Range("L98").Value = 26897
$ws.Range("N98").Value = -31389
$ws.Range("H99").Value = 12148.637
$ws.Range("I99").Value = 21804.5
$ws.Range("J99").Value = 8527.6875
$ws.Range("K99").Value = 21804.5
$ws.Range("L99").Value = 8527.6875
$ws.Range("M99").Value = -20306.5
$ws.Range("N99").Value = -11523.6875
$ws.Range("H126").Value = 12148.637
$ws.Range("I126").Value = 21804.5
$ws.Range("J126").Value = 8527.6875
$ws.Range("K126").Value = 65413.5
$ws.Range("L126").Value = 25583.0625
$ws.Range("M126").Value = -62943.5
$ws.Range("N126").Value = -30523.0625
$ws.Range("H134").Value = 2948.5
$ws.Range("I134").Value = 2235.375
$ws.Range("J134").Value = 4374.75
$ws.Range("K134").Value = 6706.125
$ws.Range("L134").Value = 13124.25
$ws.Range("M134").Value = -4171.125
$ws.Range("N134").Value = -18194.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 300
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("H32").Value = 781
$ws.Range("J32").Value = 511
$ws.Range("L32").Value = 1533
$ws.Range("N32").Value = -2099
$ws.Range("H107").Value = 2052.2856
$ws.Range("J107").Value = 953.2
$ws.Range("L107").Value = 2859.6
$ws.Range("N107").Value = -6699.6
$ws.Range("H132").Value = 4916.3335
$ws.Range("I132").Value = 666.44446
$ws.Range("J132").Value = 17666
$ws.Range("K132").Value = 5998.00014
$ws.Range("L132").Value = 158994
$ws.Range("M132").Value = -3468.00014
$ws.Range("N132").Value = -164054
$ws.Range("H136").Value = 7827.909
$ws.Range("I136").Value = 1789.6666
$ws.Range("J136").Value = 35000
$ws.Range("K136").Value = 5368.9998
$ws.Range("L136").Value = 105000
$ws.Range("M136").Value = -268.9997999999996
$ws.Range("N136").Value = -115200
$ws.Range("M9").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 345.85
$ws.Range("I2").Value = 128.18182
$ws.Range("K2").Value = 128.18182
$ws.Range("M2").Value = -15.18181999999999
$ws.Range("H80").Value = 1330412.8
$ws.Range("I80").Value = 3128997.5
$ws.Range("K80").Value = 3128997.5
$ws.Range("M80").Value = -3127999.5
$ws.Range("H83").Value = 1330412.8
$ws.Range("I83").Value = 3128997.5
$ws.Range("K83").Value = 15644987.5
$ws.Range("M83").Value = -15639995.5
$ws.Range("H97").Value = 489.7647
$ws.Range("I97").Value = 458.25925
$ws.Range("J97").Value = 611.2857
$ws.Range("K97").Value = 458.25925
$ws.Range("L97").Value = 611.2857
$ws.Range("M97").Value = 37.74074999999999
$ws.Range("N97").Value = -1603.2857
$ws.Range("H102").Value = 25010202
$ws.Range("I102").Value = 38472360
$ws.Range("J102").Value = 9049.143
$ws.Range("K102").Value = 38472360
$ws.Range("L102").Value = 9049.143
$ws.Range("M102").Value = -38470738
$ws.Range("N102").Value = -12293.143
$ws.Range("H113").Value = 1819.4546
$ws.Range("I113").Value = 1424
$ws.Range("K113").Value = 1424
$ws.Range("M113").Value = 746
$ws.Range("H126").Value = 6270.4546
$ws.Range("I126").Value = 4961.1665
$ws.Range("K126").Value = 14883.4995
$ws.Range("M126").Value = -12413.4995
$ws.Range("H132").Value = 140373.8
$ws.Range("I132").Value = 227981.22
$ws.Range("J132").Value = 8962.666999999999
$ws.Range("K132").Value = 683943.66
$ws.Range("L132").Value = 26888.001
$ws.Range("M132").Value = -681413.66
$ws.Range("N132").Value = -31948.001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 40821188
$ws.Range("I122").Value = 58827370
$ws.Range("J122").Value = 12993457
$ws.Range("K122").Value = 176482110
$ws.Range("L122").Value = 38980371
$ws.Range("M122").Value = -176479660
$ws.Range("N122").Value = -38985271
$ws.Range("H132").Value = 3295.93
$ws.Range("I132").Value = 2483.9556
$ws.Range("K132").Value = 7451.8668
$ws.Range("M132").Value = -4921.8668
$ws.Range("H136").Value = 4099.067
$ws.Range("I136").Value = 3007.9092
$ws.Range("K136").Value = 9023.7276
$ws.Range("M136").Value = -6473.7276

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 100429
$ws.Range("J46").Value = 100429
$ws.Range("L46").Value = 100429
$ws.Range("N46").Value = -100891
$ws.Range("H101").Value = 41245.75
$ws.Range("J101").Value = 41245.75
$ws.Range("L101").Value = 41245.75
$ws.Range("N101").Value = -47735.75
$ws.Range("H107").Value = 3015.5908
$ws.Range("I107").Value = 3131.4285
$ws.Range("K107").Value = 9394.2855
$ws.Range("M107").Value = -7474.2855
$ws.Range("H132").Value = 42744230
$ws.Range("I132").Value = 9260674
$ws.Range("J132").Value = 71444430
$ws.Range("K132").Value = 27782022
$ws.Range("L132").Value = 214333290
$ws.Range("M132").Value = -27779492
$ws.Range("N132").Value = -214338350
$ws.Range("H134").Value = 100429
$ws.Range("J134").Value = 100429
$ws.Range("L134").Value = 301287
$ws.Range("N134").Value = -306357
$ws.Range("H136").Value = 8445.35
$ws.Range("I136").Value = 1866.2273
$ws.Range("K136").Value = 5598.6819
$ws.Range("M136").Value = -3048.6819
